# Commit: "Can we ignore I totally missed I wrote fusionfall wrong?"
# Fixes a typo in the project title ("FussionFall" -> "FusionFall") and
# normalizes a handful of cells that were carrying a duplicate (but visually
# identical) cell style onto the shared "centered + wrap text" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the typo: "FussionFall Re:Spawn" -> "FusionFall Re:Spawn" ---
# This is the "title" column (C) of the FusionFall Re:Spawn project row.
$ws.Range("C4").Value = "FusionFall Re:Spawn"

# --- Re-apply the de-duplicated style (center/center + wrap text, default
# font) to the cells that previously referenced the redundant style index ---
$xlCenter = -4108

$cellsToRestyle = @("A3", "C3", "A4", "B4", "D4", "E4", "F4")
foreach ($addr in $cellsToRestyle) {
    $cell = $ws.Range($addr)
    $cell.HorizontalAlignment = $xlCenter
    $cell.VerticalAlignment = $xlCenter
    $cell.WrapText = $true
}
